$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.943.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.733.42'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.35'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.67%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.732.03'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.26%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.23%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.447'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000257'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.90'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.361.79'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.738.51'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.904.24'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.84'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.24%  '
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.112'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.64%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.98'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.62'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '463.81'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.693'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.59'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000146'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +8.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.16'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.23'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '29.69'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.14'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.07'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.688.06'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.138'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.993'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.74'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.88'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +15.94%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.298'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.71'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.29%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.39'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.22%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '143.90'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.53%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '385.67'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.752.99'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.90%  '
